$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.713.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.30%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.251.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.72%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'248.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.57%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.632"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.46%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'70.35"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +5.80%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.664"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +17.82%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'39.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +10.13%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'59.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.49%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0963"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.31%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'7.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +7.93%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.34%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.580.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.56%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'14.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.86%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.877"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.19%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.268.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.99%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'42.637.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.35%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0986"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.97%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.07%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'72.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.93%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'234.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.14%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.18%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +5.48%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'11.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.05%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.10%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.46%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.74%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.56%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'166.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.74%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'20.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.58%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +14.09%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +5.78%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'Hedera"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'0.0791"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +7.64%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'InjectiveProtocol"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'31.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +21.88%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +3.88%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'4.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +9.04%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'4.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.69%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +6.66%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +6.22%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'12.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +7.30%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +5.90%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'62.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.66%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'9.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +6.69%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +5.23%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'4.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.99%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +3.36%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.76%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.17%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +3.96%  "
$ws.Range("E51").Style = "Normal"
